$d = $word.ActiveDocument

$d.Content.Find.Execute("3D XPoint Apache Pass persistent memory DIMMs", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3D XPoint persistent memory DIMMs", 2)
